$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Max Christie", "SG,SF", "Dallas Mavericks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Jimmy Butler", "SF,PF", "Golden State Warriors"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Bol Bol", "PF,C", "Phoenix Suns"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Brandon Ingram", "SG,SF,PF", "Toronto Raptors"),
    @("Devin Booker", "PG,SG", "Phoenix Suns")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

Write-Output "Done updating roster."